# feat: add 2022-Q1 data
#
# 1. Insert a new worksheet "2022-Q1" between "2021-Q1" and "总计" with the
#    quarter's fund-holding detail rows.
# 2. Insert a new summary row at the top of "总计" for the 2022-Q1 totals,
#    pushing the existing 2021-Q1 summary row down.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Step 1: add the "2022-Q1" worksheet and place it right after "2021-Q1"
# (i.e. right before "总计").
# ---------------------------------------------------------------------------
$newSheet = $wb.Worksheets.Add()
$newSheet.Name = "2022-Q1"

# Re-fetch the "总计" sheet fresh (worksheet references captured before a
# collection-changing call like Add() do not track the renumbered
# collection), then move the new sheet to sit right before it.
$totalSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$newSheet.Move($totalSheet)

# Re-fetch everything by position now that the final sheet order is set:
# 1 = 2021-Q1, 2 = 2022-Q1, 3 = 总计
$ws2022 = $wb.Worksheets.Item(2)
$wsTotal = $wb.Worksheets.Item(3)

# ---------------------------------------------------------------------------
# Step 2: header row + one fund-holding row for "2022-Q1"
# ---------------------------------------------------------------------------
$ws2022.Range("B1").Value = "基金代码"
$ws2022.Range("C1").Value = "基金名称"
$ws2022.Range("D1").Value = "基金规模"
$ws2022.Range("E1").Value = "股票总仓位"
$ws2022.Range("F1").Value = "仓位占比"
$ws2022.Range("G1").Value = "持有市值(亿元)"
$ws2022.Range("H1").Value = "仓位排名"

# Match the bold / bordered / centered look used for header rows elsewhere
# in this workbook.
$header = $ws2022.Range("B1:H1")
$header.Font.Bold = $true
$header.HorizontalAlignment = -4108
$header.VerticalAlignment = -4160
$header.Borders.LineStyle = 1

# Row index cell (A2) uses the same emphasised style as the header.
$ws2022.Range("A2").Value = 0
$indexCell = $ws2022.Range("A2")
$indexCell.Font.Bold = $true
$indexCell.HorizontalAlignment = -4108
$indexCell.VerticalAlignment = -4160
$indexCell.Borders.LineStyle = 1

# Fund code / name / scale / position figures are stored as text, so force
# text formatting before writing the numeric-looking strings, then drop the
# number-format override again (the cells carry no special formatting in
# the source data, only the header/index cells do).
$ws2022.Range("B2:G2").NumberFormat = "@"
$ws2022.Range("B2").Value = "519029"
$ws2022.Range("C2").Value = "华夏稳增混合"
$ws2022.Range("D2").Value = "10.92"
$ws2022.Range("E2").Value = "92.99"
$ws2022.Range("F2").Value = "3.97"
$ws2022.Range("G2").Value = "0.4335"
$ws2022.Range("B2:G2").ClearFormats()
$ws2022.Range("H2").Value = 9

# ---------------------------------------------------------------------------
# Step 3: insert the 2022-Q1 summary row at the top of "总计", pushing the
# existing 2021-Q1 row down.
# ---------------------------------------------------------------------------
$wsTotal.Rows("2:2").Insert()

$wsTotal.Range("A2").Value = 0
$wsTotal.Range("B2").Value = "2022-Q1"
$wsTotal.Range("C2").Value = 1
$wsTotal.Range("D2").Value = 0.43

# Insert() on the row copies formatting down into the new row's cells; clear
# that so only the index cell keeps the emphasised style (matching the rest
# of the sheet, where only column A carries it).
$wsTotal.Range("B2:D2").ClearFormats()
$a2 = $wsTotal.Range("A2")
$a2.Font.Bold = $true
$a2.HorizontalAlignment = -4108
$a2.VerticalAlignment = -4160
$a2.Borders.LineStyle = 1

# The former row 2 (2021-Q1) is now row 3; its index value moves from 0 to 1.
$wsTotal.Range("A3").Value = 1
